# Update "want to go" (F column) counts across the sheets to reflect the
# latest scraped numbers, per the site's regenerated gh-pages output.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 20
$ws1.Range("F3").Value = 2740
$ws1.Range("F4").Value = 1075
$ws1.Range("F5").Value = 19858
$ws1.Range("F6").Value = 80
$ws1.Range("F7").Value = 2298
$ws1.Range("F8").Value = 754
$ws1.Range("F15").Value = 376
$ws1.Range("F16").Value = 80
$ws1.Range("F17").Value = 267
$ws1.Range("F19").Value = 199

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F14").Value = 89
$ws2.Range("F16").Value = 86

# Sheet "本地生活" (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 6010

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6010
$ws4.Range("F6").Value = 20
$ws4.Range("F8").Value = 2740
$ws4.Range("F9").Value = 1075
$ws4.Range("F10").Value = 19858
$ws4.Range("F13").Value = 80
$ws4.Range("F16").Value = 2298
$ws4.Range("F17").Value = 754
$ws4.Range("F28").Value = 376
$ws4.Range("F29").Value = 80
$ws4.Range("F32").Value = 267
$ws4.Range("F33").Value = 89
$ws4.Range("F36").Value = 199
$ws4.Range("F37").Value = 86
$ws4.Range("F38").Value = 86

$wb.Save()
